$p = $ppt.ActivePresentation

# Locate the paragraph containing the misspelled phrase "Sparse data sctructures"
# (slide 1, "Subtitle 2" placeholder) and fix it to "Sparse data structures",
# splitting the text into three separate runs: "Sparse ", "data ", "structures".
$target = "Sparse data sctructures"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange

        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi)
            if ($para.Text -notlike "*$target*") { continue }

            $offset = $para.Text.IndexOf($target)

            # "Sparse " -> characters 1..7 (left untouched, stays its own run)
            # "data "   -> characters 8..12 (retyped so it becomes its own run)
            # "sctructures" -> characters 13..23, corrected to "structures"
            $dataRun = $para.Characters($offset + 8, 5)
            $dataRun.Text = "data "

            $structRun = $para.Characters($offset + 13, 11)
            $structRun.Text = "structures"
        }
    }
}
